$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 387.5
$ws.Range("I8").Value = 387.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1162.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1023.5
$ws.Range("N8").ClearContents()
$ws.Range("H9").Value = 175.14285
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 175.14285
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 175.14285
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -513.14285
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 26316700
$ws.Range("I11").Value = 26316700
$ws.Range("K11").Value = 26316700
$ws.Range("M11").Value = -26316560
$ws.Range("H32").Value = 898.7
$ws.Range("I32").Value = 699.6667
$ws.Range("J32").Value = 984
$ws.Range("K32").Value = 699.6667
$ws.Range("L32").Value = 984
$ws.Range("M32").Value = -373.6667
$ws.Range("N32").Value = -1636
$ws.Range("H70").Value = 19235912
$ws.Range("I70").Value = 3400
$ws.Range("J70").Value = 22732732
$ws.Range("K70").Value = 10200
$ws.Range("L70").Value = 68198196
$ws.Range("M70").Value = -9930
$ws.Range("N70").Value = -68198736
$ws.Range("H73").Value = 19235912
$ws.Range("I73").Value = 3400
$ws.Range("J73").Value = 22732732
$ws.Range("K73").Value = 10200
$ws.Range("L73").Value = 68198196
$ws.Range("M73").Value = -9264
$ws.Range("N73").Value = -68200068
$ws.Range("H100").Value = 2456.9524
$ws.Range("I100").Value = 1468.4615
$ws.Range("K100").Value = 1468.4615
$ws.Range("M100").Value = -927.4614999999999
$ws.Range("H135").Value = 2261.5
$ws.Range("I135").Value = 1299
$ws.Range("K135").Value = 11691
$ws.Range("M135").Value = -9156
$ws.Range("H138").Value = 2568.0984
$ws.Range("I138").Value = 1761.2858
$ws.Range("J138").Value = 3252.6667
$ws.Range("K138").Value = 5283.857400000001
$ws.Range("L138").Value = 9758.000100000001
$ws.Range("M138").Value = -143.8574000000008
$ws.Range("N138").Value = -20038.0001
$ws.Range("H140").Value = 158776.89
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 158776.89
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 158776.89
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -169136.89
$ws.Range("H141").Value = 4124.875
$ws.Range("I141").Value = 2571.2856
$ws.Range("K141").Value = 7713.8568
$ws.Range("M141").Value = -2533.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 745
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1350
$ws.Range("H32").Value = 5869.294
$ws.Range("I32").Value = 1632.9565
$ws.Range("K32").Value = 1632.9565
$ws.Range("M32").Value = -1345.9565
$ws.Range("H45").Value = 5369.778
$ws.Range("I45").Value = 4332.7144
$ws.Range("K45").Value = 4332.7144
$ws.Range("M45").Value = -3955.7144
$ws.Range("H74").Value = 3175.3333
$ws.Range("I74").Value = 2510.9167
$ws.Range("K74").Value = 2510.9167
$ws.Range("M74").Value = -1636.9167
$ws.Range("H77").Value = 3175.3333
$ws.Range("I77").Value = 2510.9167
$ws.Range("K77").Value = 12554.5835
$ws.Range("M77").Value = -8186.583500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2842.353
$ws.Range("I20").Value = 2586.625
$ws.Range("J20").Value = 3069.6667
$ws.Range("K20").Value = 2586.625
$ws.Range("L20").Value = 3069.6667
$ws.Range("M20").Value = -2339.625
$ws.Range("N20").Value = -3563.6667
$ws.Range("H54").Value = 14800
$ws.Range("I54").Value = 14800
$ws.Range("K54").Value = 14800
$ws.Range("M54").Value = -14316
$ws.Range("H99").Value = 2375.2856
$ws.Range("I99").Value = 2362.2104
$ws.Range("K99").Value = 2362.2104
$ws.Range("M99").Value = -864.2103999999999
$ws.Range("H107").Value = 1611
$ws.Range("I107").Value = 1611
$ws.Range("K107").Value = 1611
$ws.Range("M107").Value = 309
$ws.Range("H126").Value = 118000
$ws.Range("J126").Value = 118000
$ws.Range("L126").Value = 118000
$ws.Range("N126").Value = -127880
$ws.Range("H132").Value = 56333
$ws.Range("J132").Value = 56333
$ws.Range("L132").Value = 56333
$ws.Range("N132").Value = -66453
$ws.Range("H134").Value = 13892101
$ws.Range("I134").Value = 2472.5334
$ws.Range("J134").Value = 37041480
$ws.Range("K134").Value = 7417.600199999999
$ws.Range("L134").Value = 111124440
$ws.Range("M134").Value = -4882.600199999999
$ws.Range("N134").Value = -111129510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4091.7693
$ws.Range("I31").Value = 2134.4546
$ws.Range("K31").Value = 2134.4546
$ws.Range("M31").Value = -1839.4546
$ws.Range("H34").Value = 4091.7693
$ws.Range("I34").Value = 2134.4546
$ws.Range("K34").Value = 2134.4546
$ws.Range("M34").Value = -1932.4546
$ws.Range("H41").Value = 3225.9
$ws.Range("I41").Value = 3225.9
$ws.Range("K41").Value = 3225.9
$ws.Range("M41").Value = -2797.9
$ws.Range("H58").Value = 3280.6924
$ws.Range("I58").Value = 3031.2856
$ws.Range("K58").Value = 3031.2856
$ws.Range("M58").Value = -2828.2856
$ws.Range("H86").Value = 8507
$ws.Range("I86").Value = 8507
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 8507
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -7384
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 8507
$ws.Range("I89").Value = 8507
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 42535
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -36919
$ws.Range("N89").ClearContents()
$ws.Range("H116").Value = 98992.336
$ws.Range("J116").Value = 98992.336
$ws.Range("L116").Value = 98992.336
$ws.Range("N116").Value = -108170.336
$ws.Range("H121").Value = 49000
$ws.Range("J121").Value = 49000
$ws.Range("L121").Value = 49000
$ws.Range("N121").Value = -51620
$ws.Range("H136").Value = 3280.6924
$ws.Range("I136").Value = 3031.2856
$ws.Range("K136").Value = 9093.856800000001
$ws.Range("M136").Value = -6543.856800000001
$ws.Range("H140").Value = 287799.8
$ws.Range("J140").Value = 287799.8
$ws.Range("L140").Value = 287799.8
$ws.Range("N140").Value = -298159.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 445
$ws.Range("I103").Value = 445
$ws.Range("K103").Value = 1335
$ws.Range("M103").Value = -456
$ws.Range("H131").Value = 1462.7
$ws.Range("I131").Value = 1125.3334
$ws.Range("K131").Value = 3376.0002
$ws.Range("M131").Value = 1663.9998
$ws.Range("H134").Value = 3053.5557
$ws.Range("I134").Value = 2810.25
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8430.75
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -3360.75
$ws.Range("N134").Value = -25140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 429
$ws.Range("I2").Value = 196
$ws.Range("J2").Value = 895
$ws.Range("K2").Value = 196
$ws.Range("L2").Value = 895
$ws.Range("M2").Value = -83
$ws.Range("N2").Value = -1121
$ws.Range("H15").Value = 44000
$ws.Range("J15").Value = 44000
$ws.Range("L15").Value = 44000
$ws.Range("N15").Value = -44576
$ws.Range("H81").Value = 44000
$ws.Range("J81").Value = 44000
$ws.Range("L81").Value = 44000
$ws.Range("N81").Value = -45996
$ws.Range("H84").Value = 44000
$ws.Range("J84").Value = 44000
$ws.Range("L84").Value = 132000
$ws.Range("N84").Value = -141984
$ws.Range("H99").Value = 48156.668
$ws.Range("I99").Value = 19735.5
$ws.Range("K99").Value = 19735.5
$ws.Range("M99").Value = -17489.5
$ws.Range("H114").Value = 82494.5
$ws.Range("J114").Value = 58992.668
$ws.Range("L114").Value = 58992.668
$ws.Range("N114").Value = -67670.66800000001
$ws.Range("H132").Value = 3794.2727
$ws.Range("I132").Value = 3855.7144
$ws.Range("J132").Value = 3686.75
$ws.Range("K132").Value = 11567.1432
$ws.Range("L132").Value = 11060.25
$ws.Range("M132").Value = -9037.143199999999
$ws.Range("N132").Value = -16120.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 479681.84
$ws.Range("I132").Value = 558645.5
$ws.Range("K132").Value = 1675936.5
$ws.Range("M132").Value = -1673406.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4767.75
$ws.Range("J62").Value = 4821
$ws.Range("L62").Value = 4821
$ws.Range("N62").Value = -6069
$ws.Range("H65").Value = 4767.75
$ws.Range("J65").Value = 4821
$ws.Range("L65").Value = 24105
$ws.Range("N65").Value = -30345
$ws.Range("H100").Value = 2231.2
$ws.Range("I100").Value = 2530.4285
$ws.Range("J100").Value = 1533
$ws.Range("K100").Value = 5060.857
$ws.Range("L100").Value = 3066
$ws.Range("M100").Value = -4519.857
$ws.Range("N100").Value = -4148
$ws.Range("H126").Value = 2900
$ws.Range("I126").Value = 2900
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8700
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6230
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 37885.9
$ws.Range("I132").Value = 55289.79
$ws.Range("J132").Value = 4818.5
$ws.Range("K132").Value = 165869.37
$ws.Range("L132").Value = 14455.5
$ws.Range("M132").Value = -163339.37
$ws.Range("N132").Value = -19515.5
$ws.Range("H136").Value = 27177.586
$ws.Range("I136").Value = 1941.4546
$ws.Range("K136").Value = 5824.3638
$ws.Range("M136").Value = -3274.3638
